$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 4th provider row (old row 5) entirely - the table shrinks to 3 data rows
$ws.Rows(5).Delete()

# Row 2 - Helen Kebede
$ws.Range("C2").Value = "Helen"
$ws.Range("D2").Value = "Kebede"
$ws.Range("E2").Value = "689 262 7223"
$ws.Range("F2").Value = "W854695"
$ws.Range("H2").Value = "hkebede@ottersolv.com"
$ws.Range("I2").Value = "Lesly Dorcely"

# Row 3 - Delonica James
$ws.Range("C3").Value = "Delonica"
$ws.Range("D3").Value = "James"
$ws.Range("E3").Value = "689 262 7228"
$ws.Range("F3").Value = "G183277"
$ws.Range("H3").Value = "djames@ottersolv.com"
$ws.Range("I3").Value = "Hillary Rape"

# Row 4 - Ashley Hillman
$ws.Range("C4").Value = "Ashley "
$ws.Range("D4").Value = "Hillman"
$ws.Range("E4").Value = "689 262 7229"
$ws.Range("F4").Value = "W905068"
$ws.Range("H4").Value = "ahillman@ottersolv.com"
$ws.Range("I4").Value = "Kristal Fisher"

# Rebuild the hyperlinks collection so B2:B4 keep a mailto hyperlink and the
# stale row-5 entry is gone (this engine's Hyperlinks.Delete always clears the
# whole collection, so the 3 that should remain are re-added)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:November@2024!") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:November@2024!") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:November@2024!") | Out-Null

# Re-adding hyperlinks nudges the cell style onto a freshly minted duplicate of
# the "Hyperlink" style; touching the font here snaps it back onto the
# existing shared style so B2:B4 keep their original style index
$ws.Range("B2:B4").Font.Underline = 2

# Restore the cursor/selection position recorded in the saved file
$ws.Range("F13").Select()
